# The sheet currently has headers R1..R3 (B1:D1) and R21,R22,R23,R21',R22',R23' (E1:J1).
# It needs to be expanded so E1:AC1 carry the full new header set, and the
# corresponding data columns (rows 2-6) get populated (mostly zeros, plus a
# few carried-over non-zero values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns E..AC (28 total headers span B..AC; B..D stay as-is).
$headers = @("R4","R5","R6","R7","R8","R20","R21","R22","R23","R24","R25","R26", `
             "R20X","R21X","R22X","R23X","R24X","R25X","R8X","R26X","R2X","R3X","R4X","R5X","R6X")

$startCol = 5   # column E
$styleSource = $ws.Cells.Item(1, 4)   # D1 already carries the header style
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $startCol + $i
    $cell = $ws.Cells.Item(1, $col)
    $styleSource.Copy($cell)
    $cell.Value = $headers[$i]
}

# Data rows 2-6: existing E..J values get shifted/cleared, and all of E..AC
# default to 0 except for the carried-over non-zero values below.
$lastCol = 29   # column AC
for ($row = 2; $row -le 6; $row++) {
    for ($col = $startCol; $col -le $lastCol; $col++) {
        $ws.Cells.Item($row, $col).Value = 0
    }
}

# Carried-over non-zero values (old R22 -> new R22, old R21'/R22' -> new R21X/R22X).
$ws.Cells.Item(2, 12).Value = 7    # L2 = R22
$ws.Cells.Item(2, 18).Value = 10   # R2 = R21X
$ws.Cells.Item(2, 19).Value = 5    # S2 = R22X
$ws.Cells.Item(4, 12).Value = 3    # L4 = R22
